$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "248.05"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.80"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.509"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05648"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.377"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.445"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8018"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.042"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0005824"
$ws.Range("E10").Value = "9OneONE"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1430"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07248"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03102"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.02935"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09288"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001658"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "MCDex"
$ws.Range("C17").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.216"
$ws.Range("E17").Value = "16MCDexMCB"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04719"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006501"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001050"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0003202"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.136"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.089"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04097"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006894"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1043"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009163"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005831"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7857"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.01693"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
